$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grid-search / rate-limited re-scrape re-sorted the gym listings (rows 6-21).
# Stale entries (Mercedes Club, Manhattan Plaza Health Club, CompleteBody Financial
# District) were replaced by freshly discovered ones (The Strength Club NYC, New York
# Sports Club, Gym U NYC) and the whole block was re-ordered by the new file sort.
# Column D (Email Address) is left untouched -- it is blank for every row already.
$data = @(
    ,@("Powerhouse Gym Long Island City", "https://www.powerhousegymlic.net/", "(347) 937-3030", 4.7, "40.7511444, -73.9335708", "30-30 Northern Blvd, Long Island City, NY 11101, USA")
    ,@("Life Time One Wall Street", "https://www.lifetime.life/locations/ny/one-wall-street.html?utm_source=google&utm_medium=local&utm_campaign=google-local-club", "(212) 671-7100", 4.3, "40.7067803, -74.0120018", "29 New St, New York, NY 10005, USA")
    ,@("Crunch Fitness - 34th Street", "https://www.crunch.com/locations/e-34th-st", "(212) 545-9757", 3.4, "40.7449717, -73.9770424", "222 E 34th St, New York, NY 10016, USA")
    ,@("DOGPOUND Gym, Personal Training & Group Fitness Classes", "http://thedogpound.com/", "(646) 620-6533", 4.8, "40.7245092, -74.0087489", "1 Renwick St, New York, NY 10013, USA")
    ,@("Crunch Fitness - 54th Street", "https://www.crunch.com/locations/w-54th-st", "(212) 307-7760", 4, "40.7643104, -73.9835317", "250 W 54th St, New York, NY 10019, USA")
    ,@("Crunch Fitness - FiDi", "https://www.crunch.com/locations/fidi", "(212) 308-5824", 4.1, "40.7084471, -74.00976299999999", "140 Broadway, New York, NY 10005, USA")
    ,@("Mid City Gym & Tanning", "https://midcitygym.com/", "(212) 757-0850", 4.2, "40.758403699999995, -73.9915542", "345 W 42nd St, New York, NY 10036, USA")
    ,@("GYM NYC Lexington", "https://gym-nyc.com/", "(646) 308-0192", 5, "40.7447337, -73.98133829999999", "131 E 31st St, New York, NY 10016, USA")
    ,@("CompleteBody Midtown East", "https://www.completebody.com/completebody-midtown-east/", "(212) 777-7703", 4.7, "40.7593798, -73.9646707", "301 E 57th St, New York, NY 10022, USA")
    ,@("GYM NYC East 3rd St.", "http://e3gym.nyc/", "(646) 678-5005", 4.8, "40.723565199999996, -73.98483", "155 E 3rd St, New York, NY 10009, USA")
    ,@("Powerhouse Gym 158 W. 27th Street Manhattan", "https://www.powerhousegymnyc.com/", "(646) 429-9278", 5, "40.746199999999995, -73.9934323", "158 W 27th St, New York, NY 10001, USA")
    ,@("The Strength Club NYC", "https://www.thestrengthclub.com/", "", 5, "40.7446197, -73.9868379", "1 E 28th St # 5b, New York, NY 10016, USA")
    ,@("Oculus Fitness", "http://www.oculusfitness.com/", "(646) 832-4788", 4.9, "40.7097951, -74.0063052", "151 William St, New York, NY 10038, USA")
    ,@("New York Sports Club", "https://www.newyorksportsclubs.com/clubs/23rd-park", "(212) 982-4400", 4.4, "40.7401753, -73.98571729999999", "113 E 23rd St, New York, NY 10010, USA")
    ,@("Gym U NYC", "http://gymunyc.com/", "(646) 410-2858", 4.1, "40.7445928, -73.9961144", "215 W 23rd St, New York, NY 10011, USA")
    ,@("SWEAT440 Fitness FiDi", "https://sweat440.com/gyms/new-york/financial-district-10038/?utm_source=GMB&utm_medium=organic&utm_campaign=FiDi", "(888) 507-9328", 4.9, "40.708317799999996, -74.0068091", "80 John St, New York, NY 10038, USA")
)

$startRow = 6
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]   # Name
    $ws.Cells.Item($r, 2).Value = $rowVals[1]   # Website
    $ws.Cells.Item($r, 3).Value = $rowVals[2]   # Contact Number
    $ws.Cells.Item($r, 5).Value = $rowVals[3]   # Rating
    $ws.Cells.Item($r, 6).Value = $rowVals[4]   # LatLong
    $ws.Cells.Item($r, 7).Value = $rowVals[5]   # Address
}
